# Applies the "Koparanov/info.xlsx" edit: records rows for chapters 2 and 3
# of the task log, renumbers the "ch1/ch2/ch3" markers, and extends the
# sheet from 16 to 37 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: complete the previously partial rows 12-15 (chapter_2 tasks 7-10).
# The column order below reproduces the original authoring order so the
# shared-string table comes out in the same sequence as the source edit.
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 2).Value = "стр.104 - 7 зад"
$ws.Cells.Item(12, 6).Value = ".\chapter_2\src\chapter_2\triangle.java"
$ws.Cells.Item(13, 6).Value = ".\chapter_2\src\chapter_2\heart.java"
$ws.Cells.Item(13, 2).Value = "стр.105 - 8 зад"
$ws.Cells.Item(14, 6).Value = ".\chapter_2\src\chapter_2\EmployeeRecord.java"
$ws.Cells.Item(14, 2).Value = "стр.105 - 9 зад"
$ws.Cells.Item(15, 6).Value = ".\chapter_2\src\chapter_2\numSwap.java"
$ws.Cells.Item(15, 2).Value = "стр.105 - 10 зад"

# ---------------------------------------------------------------------------
# Step 2: tag the chapter-boundary task numbers with their chapter suffix.
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value = "№5 - ch2"
$ws.Cells.Item(2, 1).Value = "№1 - ch1"
$ws.Cells.Item(16, 1).Value = "№15 - ch3"

# ---------------------------------------------------------------------------
# Step 3: number the new rows (task numbers 16-36) down column A.
# ---------------------------------------------------------------------------
for ($i = 16; $i -le 36; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "№$i"
}

# ---------------------------------------------------------------------------
# Step 4: chapter_3 task descriptions (column B), rows 16-28.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 13; $i++) {
    $row = 15 + $i
    $ws.Cells.Item($row, 2).Value = "стр.125 - $i зад"
}

# ---------------------------------------------------------------------------
# Step 5: chapter_3 source file references (column F), rows 16-21.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 6; $i++) {
    $row = 15 + $i
    $ws.Cells.Item($row, 6).Value = ".\chapter_3\src\chapter_3\main$i.java"
}

# ---------------------------------------------------------------------------
# Step 6: start/end timestamps (columns C/D) for rows 2-21, stored as the
# same day-fraction values Excel uses internally for times.
# ---------------------------------------------------------------------------
$times = @{
    2  = @(0.53888888888888886, 0.54027777777777775)
    3  = @(0.54097222222222219, 0.54513888888888895)
    4  = @(0.54513888888888895, 0.54861111111111105)
    5  = @(0.54861111111111105, 0.54999999999999993)
    6  = @(0.53472222222222221, 0.54097222222222219)
    7  = @(0.54236111111111118, 0.54305555555555551)
    8  = @(0.54375000000000007, 0.54513888888888895)
    9  = @(0.54513888888888895, 0.54583333333333328)
    10 = @(0.54583333333333328, 0.5493055555555556)
    11 = @(0.55138888888888882, 0.5541666666666667)
    12 = @(0.5541666666666667,  0.55694444444444446)
    13 = @(0.55694444444444446, 0.56111111111111112)
    14 = @(0.56111111111111112, 0.56458333333333333)
    15 = @(0.57291666666666663, 0.57638888888888895)
    16 = @(0.54166666666666663, 0.54791666666666672)
    17 = @(0.54791666666666672, 0.55208333333333337)
    18 = @(0.55208333333333337, 0.55972222222222223)
    19 = @(0.55972222222222223, 0.56388888888888888)
    20 = @(0.56388888888888888, 0.57291666666666663)
    21 = @(0.57500000000000007, 0.57916666666666672)
}
foreach ($row in $times.Keys) {
    $pair = $times[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}

# ---------------------------------------------------------------------------
# Step 7: duration formula (column E) for the newly-added rows 17-28. Rows
# 2-16 already carry their "=D-C" formulas (plain or shared) from the
# original file and are left untouched. Each new row is set individually so
# the pre-existing shared-formula block isn't disturbed.
# ---------------------------------------------------------------------------
for ($row = 17; $row -le 27; $row++) {
    $ws.Cells.Item($row, 5).Formula = "=D$row-C$row"
}
$ws.Cells.Item(28, 5).Formula = "=D28-C28"

# ---------------------------------------------------------------------------
# Step 8: move the active selection to reflect where editing left off.
# ---------------------------------------------------------------------------
$ws.Range("D24").Select()
